$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 675279.2
$ws.Range("I15").Value = 675279.2
$ws.Range("K15").Value = 2025837.6
$ws.Range("M15").Value = -2025668.6
$ws.Range("H40").Value = 3497.524
$ws.Range("J40").Value = 5313.4546
$ws.Range("L40").Value = 5313.4546
$ws.Range("N40").Value = -5663.4546
$ws.Range("H80").Value = 1091.3
$ws.Range("J80").Value = 1351.6666
$ws.Range("L80").Value = 4054.9998
$ws.Range("N80").Value = -6050.9998
$ws.Range("H83").Value = 1091.3
$ws.Range("J83").Value = 1351.6666
$ws.Range("L83").Value = 12164.9994
$ws.Range("N83").Value = -22148.9994
$ws.Range("H86").Value = 1339.5
$ws.Range("I86").Value = 1006
$ws.Range("J86").Value = 1539.6
$ws.Range("K86").Value = 1006
$ws.Range("L86").Value = 1539.6
$ws.Range("M86").Value = 117
$ws.Range("N86").Value = -3785.6
$ws.Range("H89").Value = 1339.5
$ws.Range("I89").Value = 1006
$ws.Range("J89").Value = 1539.6
$ws.Range("K89").Value = 5030
$ws.Range("L89").Value = 7698
$ws.Range("M89").Value = 586
$ws.Range("N89").Value = -18930
$ws.Range("H92").Value = 1414.375
$ws.Range("I92").Value = 1302.5
$ws.Range("J92").Value = 1750
$ws.Range("K92").Value = 1302.5
$ws.Range("L92").Value = 1750
$ws.Range("M92").Value = -54.5
$ws.Range("N92").Value = -4246
$ws.Range("H106").Value = 19619
$ws.Range("I106").Value = 15213.429
$ws.Range("K106").Value = 15213.429
$ws.Range("M106").Value = -14582.429

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 525.4286
$ws.Range("I2").Value = 525.4286
$ws.Range("K2").Value = 525.4286
$ws.Range("M2").Value = -412.4286
$ws.Range("H97").Value = 2365.9092
$ws.Range("I97").Value = 1602.5
$ws.Range("K97").Value = 1602.5
$ws.Range("M97").Value = -1106.5
$ws.Range("H102").Value = 1430.1666
$ws.Range("I102").Value = 1430.1666
$ws.Range("K102").Value = 1430.1666
$ws.Range("M102").Value = 191.8334
$ws.Range("H116").Value = 525.4286
$ws.Range("I116").Value = 525.4286
$ws.Range("K116").Value = 525.4286
$ws.Range("M116").Value = 1768.5714

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 525.4286
$ws.Range("I3").Value = 525.4286
$ws.Range("K3").Value = 525.4286
$ws.Range("M3").Value = -411.4286
$ws.Range("H20").Value = 2728.4666
$ws.Range("I20").Value = 2810.3845
$ws.Range("J20").Value = 2196
$ws.Range("K20").Value = 2810.3845
$ws.Range("L20").Value = 2196
$ws.Range("M20").Value = -2563.3845
$ws.Range("N20").Value = -2690
$ws.Range("H70").Value = 250000
$ws.Range("J70").Value = 250000
$ws.Range("L70").Value = 250000
$ws.Range("N70").Value = -250586
$ws.Range("H73").Value = 250000
$ws.Range("J73").Value = 250000
$ws.Range("L73").Value = 250000
$ws.Range("N73").Value = -252028
$ws.Range("H94").Value = 1576.9048
$ws.Range("I94").Value = 2435.8333
$ws.Range("K94").Value = 2435.8333
$ws.Range("M94").Value = -1984.8333
$ws.Range("H105").Value = 2545.5454
$ws.Range("I105").Value = 2455.6667
$ws.Range("K105").Value = 2455.6667
$ws.Range("M105").Value = -708.6667000000002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4307.75
$ws.Range("I16").Value = 1410.3334
$ws.Range("J16").Value = 13000
$ws.Range("K16").Value = 1410.3334
$ws.Range("L16").Value = 13000
$ws.Range("M16").Value = -1123.3334
$ws.Range("N16").Value = -13574
$ws.Range("H31").Value = 3849.8333
$ws.Range("I31").Value = 2700
$ws.Range("K31").Value = 2700
$ws.Range("M31").Value = -2405
$ws.Range("H34").Value = 3849.8333
$ws.Range("I34").Value = 2700
$ws.Range("K34").Value = 2700
$ws.Range("M34").Value = -2498
$ws.Range("H58").Value = 2082.4546
$ws.Range("J58").Value = 2685.6667
$ws.Range("L58").Value = 2685.6667
$ws.Range("N58").Value = -3091.6667
$ws.Range("H99").Value = 5587.8887
$ws.Range("I99").Value = 5587.8887
$ws.Range("K99").Value = 5587.8887
$ws.Range("M99").Value = -4089.8887
$ws.Range("H113").Value = 4307.75
$ws.Range("I113").Value = 1410.3334
$ws.Range("J113").Value = 13000
$ws.Range("K113").Value = 1410.3334
$ws.Range("L113").Value = 13000
$ws.Range("M113").Value = 759.6666
$ws.Range("N113").Value = -17340
$ws.Range("H126").Value = 5587.8887
$ws.Range("I126").Value = 5587.8887
$ws.Range("K126").Value = 16763.6661
$ws.Range("M126").Value = -14293.6661
$ws.Range("H134").Value = 1588
$ws.Range("I134").Value = 1384.8572
$ws.Range("J134").Value = 2299
$ws.Range("K134").Value = 4154.571599999999
$ws.Range("L134").Value = 6897
$ws.Range("M134").Value = -1619.571599999999
$ws.Range("N134").Value = -11967
$ws.Range("H136").Value = 2082.4546
$ws.Range("J136").Value = 2685.6667
$ws.Range("L136").Value = 8057.000100000001
$ws.Range("N136").Value = -13157.0001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1472.5
$ws.Range("J32").Value = 950
$ws.Range("L32").Value = 2850
$ws.Range("N32").Value = -3416
$ws.Range("H34").Value = 4655.6924
$ws.Range("J34").Value = 5912.4
$ws.Range("L34").Value = 17737.2
$ws.Range("N34").Value = -17905.2
$ws.Range("H68").Value = 900
$ws.Range("J68").Value = 900
$ws.Range("L68").Value = 2700
$ws.Range("N68").Value = -4322
$ws.Range("H71").Value = 900
$ws.Range("J71").Value = 900
$ws.Range("L71").Value = 8100
$ws.Range("N71").Value = -16212
$ws.Range("H92").Value = 274.26666
$ws.Range("I92").Value = 331
$ws.Range("J92").Value = 47.333332
$ws.Range("K92").Value = 993
$ws.Range("L92").Value = 141.999996
$ws.Range("M92").Value = 255
$ws.Range("N92").Value = -2637.999996
$ws.Range("H97").Value = 1799.8
$ws.Range("I97").Value = 1799.8
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 5399.4
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -4903.4
$ws.Range("H107").Value = 154.28572
$ws.Range("J107").Value = 154.28572
$ws.Range("L107").Value = 462.85716
$ws.Range("N107").Value = -4302.85716

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2237.6
$ws.Range("I80").Value = 2376
$ws.Range("J80").Value = 1914.6666
$ws.Range("K80").Value = 2376
$ws.Range("L80").Value = 1914.6666
$ws.Range("M80").Value = -1378
$ws.Range("N80").Value = -3910.6666
$ws.Range("H83").Value = 2237.6
$ws.Range("I83").Value = 2376
$ws.Range("J83").Value = 1914.6666
$ws.Range("K83").Value = 11880
$ws.Range("L83").Value = 9573.333000000001
$ws.Range("M83").Value = -6888
$ws.Range("N83").Value = -19557.333
$ws.Range("H97").Value = 1512.1052
$ws.Range("I97").Value = 1425.2941
$ws.Range("K97").Value = 1425.2941
$ws.Range("M97").Value = -929.2941000000001
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 2300
$ws.Range("K113").Value = 2300
$ws.Range("M113").Value = -130

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1998
$ws.Range("I7").Value = 1998
$ws.Range("K7").Value = 1998
$ws.Range("M7").Value = -1886
$ws.Range("H16").Value = 1054.3334
$ws.Range("I16").Value = 355.57144
$ws.Range("K16").Value = 355.57144
$ws.Range("M16").Value = -185.57144
$ws.Range("H61").Value = 1456.2858
$ws.Range("I61").Value = 1456.2858
$ws.Range("K61").Value = 1456.2858
$ws.Range("M61").Value = -1254.2858
$ws.Range("H113").Value = 1456.2858
$ws.Range("I113").Value = 1456.2858
$ws.Range("K113").Value = 1456.2858
$ws.Range("M113").Value = 713.7141999999999
$ws.Range("H126").Value = 1998
$ws.Range("I126").Value = 1998
$ws.Range("K126").Value = 5994
$ws.Range("M126").Value = -3524
$ws.Range("H136").Value = 15499.5
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 14177.823
$ws.Range("I136").Value = 8881.6
$ws.Range("K136").Value = 26644.8
$ws.Range("M136").Value = -24094.8
